# Update the workbook for "ESCUADRA UNION DISMAY" list:
#  - Bump the date in A1 by one month (45406 -> 45436)
#  - Update the unit prices in column D for rows 33-38

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 holds a date (serial 45406 -> 45436), keep it as a real date value
$ws.Range("A1").Value = [DateTime]::FromOADate(45436)

# Updated prices
$ws.Range("D33").Value = 203.074
$ws.Range("D34").Value = 162.027
$ws.Range("D35").Value = 151.226
$ws.Range("D36").Value = 347.818
$ws.Range("D37").Value = 248.442
$ws.Range("D38").Value = 218.196
